# Update calibration values ("Legs Update Sesi 1") on Sheet1

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 4 (L1) ---
$ws.Range("C4").Value = 1480
$ws.Range("D4").Value = 1640

# --- Row 5 (L2) ---
$ws.Range("C5").Value = 1600
$ws.Range("D5").Value = 1650
$ws.Range("E5").Value = 880
$ws.Range("F5").Value = 2100
$ws.Range("G5").Value = 2140

# --- Row 6 (L3) ---
$ws.Range("C6").Value = 1280
$ws.Range("G6").Value = 2180

# --- Row 7 (R1) ---
$ws.Range("D7").Value = 1350

# --- Row 8 (R2) ---
$ws.Range("C8").Value = 1540
$ws.Range("E8").Value = 1950

# --- Row 9 (R3) ---
$ws.Range("E9").Value = 1770
$ws.Range("F9").Value = 980
$ws.Range("G9").Value = 900

# --- Update active sheet view (scroll position / selection) ---
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 5
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("J16").Select()
